$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "2021年" row (row 13) with the same data layout as the
# preceding rows.
$ws.Cells.Item(13, 1).Value = "2021年"
$ws.Cells.Item(13, 2).Value = 20846
$ws.Cells.Item(13, 3).Value = 9724
$ws.Cells.Item(13, 4).Value = 18012.9
$ws.Cells.Item(13, 5).Value = 6118.68
$ws.Cells.Item(13, 6).Value = 15967
$ws.Cells.Item(13, 7).Value = 30.2757
$ws.Cells.Item(13, 8).Value = 32354.53
$ws.Cells.Item(13, 9).Value = 385
$ws.Cells.Item(13, 10).Value = 21179.61
$ws.Cells.Item(13, 11).Value = 29607
$ws.Cells.Item(13, 12).Value = 3941.6092
$ws.Cells.Item(13, 13).Value = 19962
$ws.Cells.Item(13, 14).Value = 53422
$ws.Cells.Item(13, 15).Value = 35.6273

# Match the formatting of column A used by the other year rows
# (bold, bordered, centered) by copying the formats from A12 onto A13.
$ws.Cells.Item(12, 1).Copy()
$ws.Cells.Item(13, 1).PasteSpecial(-4122)
$excel.CutCopyMode = 0
